$d = $word.ActiveDocument

# Table 3 (1-based index) is the "Property Definition" table with columns:
# Property name | Value type | Mandatory | Access mode | Description
$tbl = $d.Tables.Item(3)

for ($r = 2; $r -le $tbl.Rows.Count; $r++) {
    $nameCell = $tbl.Cell($r, 1).Range.Text
    if ($nameCell -match "^(value|value2)") {
        $cellRange = $tbl.Cell($r, 4).Range
        $cellRange.Find.Execute("Read Write", $false, $false, $false, $false, $false, $true, 0, $false, "", 1)
    }
}
